# Automatische test-sync: 2025-08-26 20:01:50
# Append the new "Retour status" log row (row 4) to the Logs sheet,
# extend the affected conditional-formatting ranges down to row 4,
# and bump the Dashboard category count to match.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")

# New row of data (mirrors rows 2/3 - columns C and E stay empty).
$logs.Range("A4").Value = "Retour status"
$logs.Range("B4").Value = "mailmind.test@zohomail.eu"
$logs.Range("D4").Value = "Klantenservice / Opvolging"
$logs.Range("F4").Value = "2025-08-26 20:00:54"
$logs.Range("G4").Value = "Nee"
$logs.Range("H4").Value = "Ja"
$logs.Range("I4").Value = "Nee"
$logs.Range("J4").Value = "Nee"

# Grow each existing conditional-formatting rule's range from row 3 to row 4.
$logs.Range("D2:D3").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D4"))
$logs.Range("G2:G3").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G4"))
$logs.Range("H2:H3").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H4"))
$logs.Range("I2:I3").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I4"))
$logs.Range("J2:J3").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J4"))

# Dashboard count for "Klantenservice / Opvolging" goes from 2 to 3.
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 3
